$d = $word.ActiveDocument

# Remove every comment in the document (and their range markers/references).
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}

# Word stamps the "_GoBack" bookmark at the location of the most recent edit.
# Re-anchor it at the point of the last substantive change in this revision
# (the start of the "There is not enough evidence..." sentence).
$r = $d.Content
$r.Find.Execute("There is not enough evidence", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(1)
$d.Bookmarks.Add("_GoBack", $r)
